$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.91"
$ws.Range("E2").Value = "'-3.41%"
$ws.Range("D3").Value = "'41.77"
$ws.Range("E3").Value = "'-5.09%"
$ws.Range("D4").Value = "'5.179"
$ws.Range("E4").Value = "'-3.93%"
$ws.Range("D5").Value = "'0.08107"
$ws.Range("E5").Value = "'-3.25%"
$ws.Range("D6").Value = "'4.366"
$ws.Range("E6").Value = "'-1.56%"
$ws.Range("D7").Value = "'1.745"
$ws.Range("E7").Value = "'-10.55%"
$ws.Range("D8").Value = "'0.9303"
$ws.Range("E8").Value = "'-4.73%"
$ws.Range("D9").Value = "'0.1124"
$ws.Range("E9").Value = "'-1.37%"
$ws.Range("D10").Value = "'0.1859"
$ws.Range("E10").Value = "'-2.30%"
$ws.Range("D11").Value = "'0.09293"
$ws.Range("E11").Value = "'-5.03%"
$ws.Range("D12").Value = "'0.04578"
$ws.Range("E12").Value = "'-1.04%"
$ws.Range("D13").Value = "'7.412"
$ws.Range("E13").Value = "'-19.09%"
$ws.Range("E14").Value = "'-0.86%"
$ws.Range("D15").Value = "'0.001281"
$ws.Range("E15").Value = "'-0.37%"
$ws.Range("D16").Value = "'0.005849"
$ws.Range("E16").Value = "'-4.34%"
$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'-1.49%"
$ws.Range("E18").Value = "'0.94%"
$ws.Range("D19").Value = "'0.3373"
$ws.Range("E19").Value = "'1.28%"
$ws.Range("D20").Value = "'0.1381"
$ws.Range("E20").Value = "'0.86%"
$ws.Range("D21").Value = "'0.2598"
$ws.Range("D22").Value = "'0.04187"
$ws.Range("E22").Value = "'0.72%"
$ws.Range("D23").Value = "'0.001245"
$ws.Range("E23").Value = "'-3.76%"
$ws.Range("D24").Value = "'0.004269"
$ws.Range("E24").Value = "'-3.25%"
$ws.Range("D25").Value = "'0.0001226"
$ws.Range("E25").Value = "'-5.82%"
$ws.Range("D26").Value = "'0.0002984"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("D38").Value = "'0.02599"
$ws.Range("E38").Value = "'-2.35%"
$ws.Range("D39").Value = "'0.05483"
$ws.Range("E39").Value = "'-2.80%"
$ws.Range("D40").Value = "'0.008061"
$ws.Range("E40").Value = "'3.21%"
$ws.Range("D41").Value = "'0.1393"
$ws.Range("E41").Value = "'-1.60%"
$ws.Range("D42").Value = "'0.006514"
$ws.Range("D43").Value = "'0.002094"
$ws.Range("E43").Value = "'-0.85%"
$ws.Range("D44").Value = "'0.008262"
$ws.Range("E44").Value = "'4.41%"
$ws.Range("D45").Value = "'0.3450"
$ws.Range("E45").Value = "'-1.81%"
$ws.Range("D46").Value = "'0.00006724"
$ws.Range("E46").Value = "'-1.34%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("D48").Value = "'0.003400"
$ws.Range("E48").Value = "'-3.24%"
$ws.Range("D49").Value = "'0.004108"
$ws.Range("E49").Value = "'16.36%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("E51").Value = "'0.08%"
